$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 64

# Insertion order matches target shared-strings order: title, link, author
$ws.Cells.Item($row, 2).Value = "Creating beautiful tables in R with {gt}"
$ws.Cells.Item($row, 3).Value = "https://gt.albert-rapp.de/"
$ws.Cells.Item($row, 1).Value = "Albert Rapp"

$ws.Range("A" + $row + ":C" + $row).WrapText = $true
$ws.Rows.Item($row).RowHeight = 34

$ws.Range("B64").Select()
$excel.ActiveWindow.ScrollRow = 60
